# Wyoming Minerals Initiative Concept Paper — highlight two recently
# updated terms ("post-mining" and "Minerals") to call out the latest
# Award-data refresh (commit: "Update with latest Award data (till 07082025)").

$d = $word.ActiveDocument

# --- Edit 1 -----------------------------------------------------------
# "...fostering development of much of the post-mining supply chain in
# Wyoming." -> highlight just the word "post-mining" in yellow.
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("much of the post-mining supply chain in Wyoming. ")
if ($found1) {
    $scope1 = $d.Range($rng1.Start, $rng1.End)
    $find1 = $scope1.Find
    $find1.ClearFormatting()
    $find1.Replacement.ClearFormatting()
    $find1.Replacement.Highlight = $true
    $find1.Execute("post-mining", $true, $false, $false, $false, $false, $true, 1, $false, "post-mining", 2)
}

# --- Edit 2 -------------------------------------------------------------
# Table cell: "Wyoming Minerals Quasi-endowment for seed grants" -> only
# highlight the word "Minerals" in yellow.
$rng2 = $d.Content
$found2 = $rng2.Find.Execute(" Minerals Quasi-endowment")
if ($found2) {
    $scope2 = $d.Range($rng2.Start, $rng2.End)
    $find2 = $scope2.Find
    $find2.ClearFormatting()
    $find2.Replacement.ClearFormatting()
    $find2.Replacement.Highlight = $true
    $find2.Execute("Minerals", $true, $false, $false, $false, $false, $true, 1, $false, "Minerals", 2)
}
